$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links) - safe to assign directly.
$textUpdates = @(
    @{ Cell = 'B8'; Value = 'GateToken' },
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'B9'; Value = 'MXToken' },
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'B11'; Value = 'WazirX' },
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' },
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'B13'; Value = 'BitrueCoin' },
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'B14'; Value = 'BitMartToken' },
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'B15'; Value = 'BitForexToken' },
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'B16'; Value = 'TigerCash' },
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'B17'; Value = 'LEO' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
)

foreach ($item in $textUpdates) {
    $ws.Range($item.Cell).Value = $item.Value
}

# Numeric-looking text updates (prices, percentages) - these must stay as
# literal text (matching the source inlineStr cells), not get auto-coerced
# into numbers/percent-formatted values by Excel. Force text format, write
# the value, then restore the Normal style so no stray formatting sticks.
$numericTextUpdates = @(
    @{ Cell = 'D2'; Value = '292.84' },
    @{ Cell = 'E2'; Value = '-0.06%' },
    @{ Cell = 'D3'; Value = '31.25' },
    @{ Cell = 'E3'; Value = '1.11%' },
    @{ Cell = 'D4'; Value = '4.958' },
    @{ Cell = 'E4'; Value = '1.05%' },
    @{ Cell = 'D5'; Value = '0.07499' },
    @{ Cell = 'E5'; Value = '2.62%' },
    @{ Cell = 'D6'; Value = '2.274' },
    @{ Cell = 'E6'; Value = '-1.01%' },
    @{ Cell = 'D7'; Value = '7.801' },
    @{ Cell = 'E7'; Value = '1.38%' },
    @{ Cell = 'D8'; Value = '3.769' },
    @{ Cell = 'E8'; Value = '1.11%' },
    @{ Cell = 'D9'; Value = '0.9202' },
    @{ Cell = 'E9'; Value = '2.14%' },
    @{ Cell = 'D10'; Value = '0.09355' },
    @{ Cell = 'E10'; Value = '18.40%' },
    @{ Cell = 'D11'; Value = '0.1738' },
    @{ Cell = 'E11'; Value = '3.40%' },
    @{ Cell = 'D12'; Value = '0.08374' },
    @{ Cell = 'E12'; Value = '3.99%' },
    @{ Cell = 'D13'; Value = '0.03275' },
    @{ Cell = 'E13'; Value = '5.56%' },
    @{ Cell = 'D14'; Value = '0.09942' },
    @{ Cell = 'E14'; Value = '-1.08%' },
    @{ Cell = 'D15'; Value = '0.001501' },
    @{ Cell = 'E15'; Value = '0.23%' },
    @{ Cell = 'D16'; Value = '0.005797' },
    @{ Cell = 'D17'; Value = '3.475' },
    @{ Cell = 'E17'; Value = '-0.07%' },
    @{ Cell = 'E19'; Value = '0.43%' },
    @{ Cell = 'E20'; Value = '1.06%' },
    @{ Cell = 'D21'; Value = '4.100' },
    @{ Cell = 'E21'; Value = '2.12%' },
    @{ Cell = 'D22'; Value = '0.2099' },
    @{ Cell = 'E22'; Value = '0.06%' },
    @{ Cell = 'D23'; Value = '0.04539' },
    @{ Cell = 'E23'; Value = '0.08%' },
    @{ Cell = 'D24'; Value = '0.001219' },
    @{ Cell = 'E24'; Value = '0.66%' },
    @{ Cell = 'D25'; Value = '0.004311' },
    @{ Cell = 'E25'; Value = '-6.96%' },
    @{ Cell = 'E26'; Value = '0.08%' },
    @{ Cell = 'D27'; Value = '0.0003391' },
    @{ Cell = 'E27'; Value = '0.02%' },
    @{ Cell = 'D39'; Value = '0.01643' },
    @{ Cell = 'E39'; Value = '2.91%' },
    @{ Cell = 'D40'; Value = '0.04592' },
    @{ Cell = 'E40'; Value = '3.59%' },
    @{ Cell = 'E41'; Value = '1.72%' },
    @{ Cell = 'D42'; Value = '0.009831' },
    @{ Cell = 'E42'; Value = '14.06%' },
    @{ Cell = 'D43'; Value = '0.1361' },
    @{ Cell = 'E43'; Value = '2.98%' },
    @{ Cell = 'D44'; Value = '0.002218' },
    @{ Cell = 'E44'; Value = '10.96%' },
    @{ Cell = 'D45'; Value = '0.009410' },
    @{ Cell = 'E45'; Value = '-0.17%' },
    @{ Cell = 'D46'; Value = '0.00006099' },
    @{ Cell = 'E46'; Value = '2.93%' },
    @{ Cell = 'E47'; Value = '0.02%' },
    @{ Cell = 'D48'; Value = '2.794' },
    @{ Cell = 'E48'; Value = '24.69%' },
    @{ Cell = 'D49'; Value = '0.001998' },
    @{ Cell = 'E49'; Value = '-30.98%' },
    @{ Cell = 'D50'; Value = '0.00002098' },
    @{ Cell = 'E50'; Value = '0.02%' },
    @{ Cell = 'D51'; Value = '0.0001998' },
    @{ Cell = 'E51'; Value = '0.02%' }
)

foreach ($item in $numericTextUpdates) {
    $c = $ws.Range($item.Cell)
    $c.NumberFormat = "@"
    $c.Value = $item.Value
    $c.Style = "Normal"
}
